$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.390.09"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "3.034.72"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'578.28"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").Value = "'168.33"
$ws.Range("E6").Value = "  +2.71%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.032.52"
$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").Value = "'0.520"
$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("E10").Value = "  +2.99%  "

$ws.Range("E11").Value = "  -1.46%  "

$ws.Range("D12").Value = "'0.484"
$ws.Range("E12").Value = "  +6.04%  "

$ws.Range("E13").Value = "  -2.09%  "

$ws.Range("D14").Value = "'36.66"
$ws.Range("E14").Value = "  +5.75%  "

$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("D16").Value = "66.335.12"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").Value = "3.532.54"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").Value = "'7.26"
$ws.Range("E18").Value = "  +4.70%  "

$ws.Range("D19").Value = "'16.46"
$ws.Range("E19").Value = "  +18.74%  "

$ws.Range("D20").Value = "3.026.75"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").Value = "'465.65"
$ws.Range("E21").Value = "  +2.36%  "

$ws.Range("D22").Value = "'0.709"
$ws.Range("E22").Value = "  +3.15%  "

$ws.Range("D23").Value = "'7.39"
$ws.Range("E23").Value = "  +0.55%  "

$ws.Range("D24").Value = "'83.01"
$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("E25").Value = "  +4.03%  "

$ws.Range("D26").Value = "'2.27"
$ws.Range("E26").Value = "  -1.91%  "

$ws.Range("D27").Value = "'10.07"
$ws.Range("E27").Value = "  -3.35%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").Value = "  +0.78%  "

$ws.Range("D30").Value = "'2.46"
$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("D31").Value = "'2.64"
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("E32").Value = "  +7.08%  "

$ws.Range("E33").Value = "  -4.48%  "

$ws.Range("D34").Value = "'28.36"
$ws.Range("E34").Value = "  +3.47%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").Value = "'5.90"
$ws.Range("E36").Value = "  +1.19%  "

$ws.Range("D37").Value = "'0.993"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").Value = "'48.95"
$ws.Range("E38").Value = "  +11.20%  "

$ws.Range("D39").Value = "'2.06"
$ws.Range("E39").Value = "  -1.67%  "

$ws.Range("E42").Value = "  -0.86%  "

$ws.Range("D43").Value = "'2.85"
$ws.Range("E43").Value = "  -4.49%  "

$ws.Range("D44").Value = "'8.62"
$ws.Range("E44").Value = "  +2.12%  "

$ws.Range("D45").Value = "'0.0361"
$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("D46").Value = "'380.63"
$ws.Range("E46").Value = "  -5.62%  "

$ws.Range("D47").Value = "2.711.61"
$ws.Range("E47").Value = "  -2.82%  "

$ws.Range("D48").Value = "'135.02"
$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("D50").Value = "'24.57"
$ws.Range("E50").Value = "  +2.80%  "

$ws.Range("D51").Value = "'2.23"
$ws.Range("E51").Value = "  +3.43%  "

# Swap rows 40 and 41 (OKB <-> TheGraph)
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.313"
$ws.Range("E40").Value = "  +0.88%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'49.47"
$ws.Range("E41").Value = "  -0.58%  "

